# Add a new "status" column (C) to the Data sheet, writing a header and
# a "Pass" value for the existing data row, mirroring the "write data
# into excel" functionality added by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "status"
$ws.Range("C2").Value = "Pass"

# Leave the selection on the newly written header cell, matching the
# post-edit workbook state.
$ws.Range("C1").Select() | Out-Null
